$d = $word.ActiveDocument

# Locate the existing "@jakeg82" text.
$find = $d.Content
$find.Find.ClearFormatting()
$find.Find.Text = "@jakeg82"
$found = $find.Find.Execute()

if ($found) {
    $atStart = $find.Start

    # Range covering just the "@" character (stays untouched, becomes the
    # first run) and the range covering "jakeg82" (to be replaced).
    $restRange = $d.Range($atStart + 1, $find.End)
    $restRange.Text = "not-used"

    # Force the run to split from the "@" run by toggling direct
    # character formatting on/off (net formatting change is a no-op, but
    # it causes the engine to materialize a separate <w:r> for the new
    # text instead of merging it back into the previous run).
    $newTextRange = $d.Range($atStart + 1, $atStart + 1 + 8)
    $newTextRange.Bold = 1
    $newTextRange.Bold = 0
}
